$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A ("ФИО") before the existing "Класс" column
$ws.Range("A1").EntireColumn.Insert()

# Insert two new columns ("Название", "Предмет") before the existing "Этап" column
# After inserting column A, former column D ("Этап") is now column E
$ws.Range("E1:F1").EntireColumn.Insert()

# Clear old data rows (previously rows 2-5, still rows 2-5 after column inserts)
$ws.Range("A2:H5").Clear()

# Header row
$ws.Range("A1").Value = "ФИО"
$ws.Range("B1").Value = "Класс"
$ws.Range("C1").Value = "Категории"
$ws.Range("D1").Value = "Тип"
$ws.Range("E1").Value = "Название"
$ws.Range("F1").Value = "Предмет"
$ws.Range("G1").Value = "Этап"
$ws.Range("H1").Value = "Результат"

# Data row 2
$ws.Range("A2").Value = "s"
$ws.Range("B2").Value = "11-А"
$ws.Range("C2").Value = "Интеллектуальные соревнования"
$ws.Range("D2").Value = "МОШ"
$ws.Range("E2").Value = "Int"
$ws.Range("F2").Value = "Int"
$ws.Range("G2").Value = "отборочный"
$ws.Range("H2").Value = "победитель"

# Data row 3
$ws.Range("A3").Value = "s"
$ws.Range("B3").Value = "11-А"
$ws.Range("C3").Value = "Спортивные достижения"
$ws.Range("D3").Value = "Турнир по баскетболу"
$ws.Range("E3").Value = "ооо"
$ws.Range("F3").Value = "моя оборона"
$ws.Range("G3").Value = "школьный"
$ws.Range("H3").Value = "победитель"

# Column widths matching target layout (values chosen so the engine's
# character-width quantization lands as close as possible to the target
# stored "width" attributes of 4.570313 / 6.998291 / 35.2771 / 24.708252 /
# 10.568848 / 13.996582 / 12.854004 / 12.854004)
$ws.Columns.Item(1).ColumnWidth = 3.666667
$ws.Columns.Item(2).ColumnWidth = 6.166667
$ws.Columns.Item(3).ColumnWidth = 34.5
$ws.Columns.Item(4).ColumnWidth = 23.833333
$ws.Columns.Item(5).ColumnWidth = 9.666667
$ws.Columns.Item(6).ColumnWidth = 13.166667
$ws.Columns.Item(7).ColumnWidth = 12.0
$ws.Columns.Item(8).ColumnWidth = 12.0
